# Responsible Teaming.pptx -- "Add files via upload" / PDD revision
# Slide 4 ("Team Working Agreement Pet4All") gets:
#   1. Title retitled to "Team Working Part of Product Description Doc (PDD)"
#   2. Heading line "PET4ALL" -> "PET4ALL - " + new bold/coloured run
#      "Requested Approach/Behavior"
#   3. First bullet split into "Only " + "one person talks at a time and we listen"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- 1. Title -------------------------------------------------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Characters(1, $titleRange.Length).Text = "Team Working Part of Product Description Doc (PDD)"

# --- 2. Body placeholder ----------------------------------------------------
$body = $s.Shapes.Item(2)
$bodyRange = $body.TextFrame.TextRange

# 2a. "PET4ALL" -> "PET4ALL - "
$headingPara = $bodyRange.Paragraphs(1)
$headingPara.Characters(1, $headingPara.Length).Text = "PET4ALL - "

# 2b. append the new "Requested Approach/Behavior" run after the heading
$headingPara = $bodyRange.Paragraphs(1)
$oldLen = $headingPara.Length
$grown = $headingPara.InsertAfter("Requested Approach/Behavior")
$newRun = $grown.Characters($oldLen, $grown.Length - $oldLen)
$newRun.Font.Bold = $true
$newRun.Font.Italic = $false
$newRun.Font.Color.RGB = 0xD9D1C9
$newRun.Font.Name = "-apple-system"

# 2c. split "Only one person talks at a time and we listen" into two runs
$bullet1 = $bodyRange.Paragraphs(2)
$prefix = $bullet1.Characters(1, 5)
$prefix.Text = "Only "
